# [WIP] Added REPCA1 model; added `service.CurrentSign`
#
# Adds two new worksheets to the ANDES kundur_reg.xlsx case file:
#   - "BusFreq" (bus frequency measurement device), inserted before "REGCA1"
#   - "REPCA1"  (renewable plant control model), inserted after "REECA1"
#     (i.e. right before "Toggler")
#
# Final sheet order:
#   Bus, PQ, PV, Slack, Line, Area, GENROU, TGOV1, EXDC2,
#   BusFreq, REGCA1, REECA1, REPCA1, Toggler

$wb = $excel.ActiveWorkbook

function ColToNum($col) {
    $num = 0
    foreach ($c in $col.ToCharArray()) {
        $num = $num * 26 + ([int][char]$c - [int][char]'A' + 1)
    }
    return $num
}

function Style-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Style-UidCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Write-ModelSheet($ws, $headers, $rowData) {
    # Header row (row 1), every cell styled like the rest of the workbook's
    # model sheets (bold, thin border, centered).
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item(1, $col)
        $cell.Value = $headers[$i]
        Style-HeaderCell $cell
    }

    # Data rows start at row 2. "uid" (column A) is always styled; the rest
    # of the columns use the default style, matching the existing sheets.
    $r = 2
    foreach ($row in $rowData) {
        foreach ($key in $row.Keys) {
            $col = ColToNum $key
            $cell = $ws.Cells.Item($r, $col)
            $cell.Value = $row[$key]
        }
        Style-UidCell ($ws.Cells.Item($r, 1))
        $r++
    }

    # Freeze the header row, as every other model sheet in the workbook does.
    $ws.Activate()
    $ws.Range("A2").Select()
    $excel.ActiveWindow.FreezePanes = $true
}

# ---------------------------------------------------------------------------
# 1. BusFreq - inserted immediately before REGCA1
# ---------------------------------------------------------------------------
$regca1 = $wb.Worksheets.Item("REGCA1")
$busFreq = $wb.Worksheets.Add($regca1)
$busFreq.Name = "BusFreq"

$busFreqHeaders = @("uid", "idx", "u", "name", "bus", "Tf", "Tw", "fn")
$busFreqRows = @(
    [ordered]@{
        "A" = 0
        "B" = "BusFreq_1"
        "C" = 1
        "D" = "BusFreq_1"
        "E" = 4
        "F" = 0.02
        "G" = 0.02
        "H" = 60
    }
)
Write-ModelSheet $busFreq $busFreqHeaders $busFreqRows
$busFreq.Range("M28").Select()

# ---------------------------------------------------------------------------
# 2. REPCA1 - inserted immediately after REECA1 (i.e. before Toggler)
# ---------------------------------------------------------------------------
$toggler = $wb.Worksheets.Item("Toggler")
$repca1 = $wb.Worksheets.Add($toggler)
$repca1.Name = "REPCA1"

$repca1Headers = @(
    "uid", "idx", "u", "name",
    "ree", "line", "busr", "busf", "Tfltr",
    "Kp", "Ki", "Tft", "Tfv", "Vfrz", "Rc", "Xc", "emax", "emin",
    "dbd1", "dbd2", "Qmax", "Qmin", "Kpg", "Kig", "Tp",
    "fdbd1", "fdbd2", "femax", "femin", "Pmax", "Pmin", "Tg",
    "Ddn", "Dup"
)
$repca1Rows = @(
    [ordered]@{
        "A" = 0
        "B" = "REPCA1_1"
        "C" = 1
        "D" = "REPCA1_1"
        "E" = 1
        "F" = "Line_14"
        "H" = "BusFreq_1"
        "I" = 0.02
        "J" = 1
        "K" = 0.1
        "L" = 1
        "M" = 1
        "N" = 0.8
        "O" = 0
        "P" = 0.01
        "Q" = 999
        "R" = -999
        "S" = -0.1
        "T" = 0.1
        "U" = 999
        "V" = -999
        "W" = 1
        "X" = 0.1
        "Y" = 0.02
        "Z" = -0.01
        "AA" = 0.01
        "AB" = 0.05
        "AC" = -0.05
        "AD" = 999
        "AE" = 0
        "AF" = 0.02
        "AG" = 0.05
        "AH" = 0.05
    }
)
Write-ModelSheet $repca1 $repca1Headers $repca1Rows

# ---------------------------------------------------------------------------
# Final view state: BusFreq was the last-active sheet when the workbook was
# saved (matches the author's activeTab / tabSelected in the target file).
# ---------------------------------------------------------------------------
$busFreq.Activate()
